# Insert a new weekly record at row 349 in "Sheet1" (the Papa/Hortaliza
# price sheet). All existing rows from 349 onward shift down by one,
# extending the used range from A1:R375 to A1:R376.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 349 (and everything below it) down by one row.
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new record.
$ws.Cells.Item(349, 1).Value = 7
$ws.Cells.Item(349, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(349, 3).Value = "Ñuble"
$ws.Cells.Item(349, 4).Value = 44783
$ws.Cells.Item(349, 5).Value = 16
$ws.Cells.Item(349, 6).Value = 100114001
$ws.Cells.Item(349, 7).Value = "Papa"
$ws.Cells.Item(349, 8).Value = "Patagonia"
$ws.Cells.Item(349, 9).Value = "1a (guarda)"
$ws.Cells.Item(349, 10).Value = 160
$ws.Cells.Item(349, 11).Value = 7000
$ws.Cells.Item(349, 12).Value = 8000
$ws.Cells.Item(349, 13).Value = 7500
$ws.Cells.Item(349, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(349, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(349, 16).Value = 300
$ws.Cells.Item(349, 17).Value = 25
$ws.Cells.Item(349, 18).Value = "Hortaliza"
